$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 90.375
$ws.Range("I4").Value = 105.5
$ws.Range("K4").Value = 105.5
$ws.Range("M4").Value = 8.5
$ws.Range("H12").Value = 11725.777
$ws.Range("J12").Value = 1798
$ws.Range("L12").Value = 1798
$ws.Range("N12").Value = -2138
$ws.Range("H15").Value = 342951.38
$ws.Range("I15").Value = 342951.38
$ws.Range("K15").Value = 1028854.14
$ws.Range("M15").Value = -1028685.14
$ws.Range("H86").Value = 21683.25
$ws.Range("I86").Value = 28816.334
$ws.Range("J86").Value = 14550.167
$ws.Range("K86").Value = 28816.334
$ws.Range("L86").Value = 14550.167
$ws.Range("M86").Value = -27693.334
$ws.Range("N86").Value = -16796.167
$ws.Range("H89").Value = 21683.25
$ws.Range("I89").Value = 28816.334
$ws.Range("J89").Value = 14550.167
$ws.Range("K89").Value = 144081.67
$ws.Range("L89").Value = 72750.83499999999
$ws.Range("M89").Value = -138465.67
$ws.Range("N89").Value = -83982.83499999999
$ws.Range("H132").Value = 4095.7273
$ws.Range("I132").Value = 4095.7273
$ws.Range("K132").Value = 12287.1819
$ws.Range("M132").Value = -9757.1819
$ws.Range("H137").Value = 5017.5
$ws.Range("I137").Value = 1681.5
$ws.Range("J137").Value = 7241.5
$ws.Range("K137").Value = 5044.5
$ws.Range("L137").Value = 21724.5
$ws.Range("M137").Value = -2494.5
$ws.Range("N137").Value = -26824.5
$ws.Range("H138").Value = 7649
$ws.Range("J138").Value = 8017.343
$ws.Range("L138").Value = 24052.029
$ws.Range("N138").Value = -34332.02899999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 84000
$ws.Range("J7").Value = 84000
$ws.Range("L7").Value = 84000
$ws.Range("N7").Value = -84228
$ws.Range("H32").Value = 35523.06
$ws.Range("I32").Value = 11648.667
$ws.Range("J32").Value = 48545.453
$ws.Range("K32").Value = 11648.667
$ws.Range("L32").Value = 48545.453
$ws.Range("M32").Value = -11361.667
$ws.Range("N32").Value = -49119.453
$ws.Range("H61").Value = 125005384
$ws.Range("I61").Value = 125005384
$ws.Range("K61").Value = 125005384
$ws.Range("M61").Value = -125005172
$ws.Range("H97").Value = 1589.9231
$ws.Range("I97").Value = 1527.875
$ws.Range("K97").Value = 1527.875
$ws.Range("M97").Value = -1031.875
$ws.Range("H136").Value = 125005384
$ws.Range("I136").Value = 125005384
$ws.Range("K136").Value = 375016152
$ws.Range("M136").Value = -375013602

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2089.5
$ws.Range("I99").Value = 1771.4286
$ws.Range("K99").Value = 1771.4286
$ws.Range("M99").Value = -273.4286
$ws.Range("H134").Value = 3855.9656
$ws.Range("I134").Value = 3707.9644
$ws.Range("K134").Value = 11123.8932
$ws.Range("M134").Value = -8588.893199999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10376
$ws.Range("I31").Value = 6821.375
$ws.Range("K31").Value = 6821.375
$ws.Range("M31").Value = -6526.375
$ws.Range("H34").Value = 10376
$ws.Range("I34").Value = 6821.375
$ws.Range("K34").Value = 6821.375
$ws.Range("M34").Value = -6619.375
$ws.Range("H58").Value = 9551
$ws.Range("I58").Value = 10279.5
$ws.Range("J58").Value = 4451.5
$ws.Range("K58").Value = 10279.5
$ws.Range("L58").Value = 4451.5
$ws.Range("M58").Value = -10076.5
$ws.Range("N58").Value = -4857.5
$ws.Range("H60").Value = 49999
$ws.Range("J60").Value = 49999
$ws.Range("L60").Value = 49999
$ws.Range("N60").Value = -51021
$ws.Range("H132").Value = 2506.3157
$ws.Range("I132").Value = 2271.9412
$ws.Range("J132").Value = 4498.5
$ws.Range("K132").Value = 6815.823600000001
$ws.Range("L132").Value = 13495.5
$ws.Range("M132").Value = -4285.823600000001
$ws.Range("N132").Value = -18555.5
$ws.Range("H136").Value = 9551
$ws.Range("I136").Value = 10279.5
$ws.Range("J136").Value = 4451.5
$ws.Range("K136").Value = 30838.5
$ws.Range("L136").Value = 13354.5
$ws.Range("M136").Value = -28288.5
$ws.Range("N136").Value = -18454.5
$ws.Range("H141").Value = 633771
$ws.Range("J141").Value = 684778.9
$ws.Range("L141").Value = 684778.9
$ws.Range("N141").Value = -695138.9

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 584.7692
$ws.Range("J12").Value = 644.2857
$ws.Range("L12").Value = 1932.8571
$ws.Range("N12").Value = -2278.8571
$ws.Range("H109").Value = 3000
$ws.Range("I109").Value = 3000
$ws.Range("K109").Value = 9000
$ws.Range("M109").Value = -7960
$ws.Range("H122").Value = 2939.8
$ws.Range("I122").Value = 1298.5
$ws.Range("J122").Value = 4034
$ws.Range("K122").Value = 11686.5
$ws.Range("L122").Value = 36306
$ws.Range("M122").Value = -9236.5
$ws.Range("N122").Value = -41206
$ws.Range("H141").Value = 3217.1177
$ws.Range("I141").Value = 2638.1538
$ws.Range("J141").Value = 5098.75
$ws.Range("K141").Value = 7914.4614
$ws.Range("L141").Value = 15296.25
$ws.Range("M141").Value = -2734.4614
$ws.Range("N141").Value = -25656.25

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2660
$ws.Range("J80").Value = 2400
$ws.Range("L80").Value = 2400
$ws.Range("N80").Value = -4396
$ws.Range("H83").Value = 2660
$ws.Range("J83").Value = 2400
$ws.Range("L83").Value = 12000
$ws.Range("N83").Value = -21984
$ws.Range("H132").Value = 8731.277
$ws.Range("I132").Value = 8727.444
$ws.Range("K132").Value = 26182.332
$ws.Range("M132").Value = -23652.332

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2746.611
$ws.Range("I16").Value = 644
$ws.Range("K16").Value = 644
$ws.Range("M16").Value = -474
$ws.Range("H82").Value = 2591.7273
$ws.Range("I82").Value = 2397
$ws.Range("J82").Value = 2932.5
$ws.Range("K82").Value = 2397
$ws.Range("L82").Value = 2932.5
$ws.Range("M82").Value = -2036
$ws.Range("N82").Value = -3654.5
$ws.Range("H85").Value = 2591.7273
$ws.Range("I85").Value = 2397
$ws.Range("J85").Value = 2932.5
$ws.Range("K85").Value = 2397
$ws.Range("L85").Value = 2932.5
$ws.Range("M85").Value = -1149
$ws.Range("N85").Value = -5428.5
$ws.Range("H102").Value = 69999
$ws.Range("J102").Value = 69999
$ws.Range("L102").Value = 69999
$ws.Range("N102").Value = -76489
$ws.Range("H132").Value = 46014660
$ws.Range("I132").Value = 51126736
$ws.Range("K132").Value = 153380208
$ws.Range("M132").Value = -153377678
$ws.Range("H136").Value = 3173.75
$ws.Range("I136").Value = 3173.75
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 9521.25
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -6971.25
$ws.Range("N136").ClearContents()  # remove cell (was -14097)

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 42376
$ws.Range("J27").Value = 42376
$ws.Range("L27").Value = 42376
$ws.Range("N27").Value = -42514
$ws.Range("H100").Value = 1562
$ws.Range("I100").Value = 1452.2
$ws.Range("K100").Value = 2904.4
$ws.Range("M100").Value = -2363.4
$ws.Range("H136").Value = 2443.7144
$ws.Range("I136").Value = 2467.4443
$ws.Range("K136").Value = 7402.3329
$ws.Range("M136").Value = -4852.3329
